$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of cell address -> new text value (values are stored as text/strings
# in the sheet, matching the original "Price" and "Volume(1h)" columns)
$updates = @{
    'D2' = '332.11'
    'E2' = '0.52%'
    'D3' = '45.50'
    'E3' = '2.65%'
    'D4' = '5.552'
    'E4' = '0.93%'
    'D5' = '0.08368'
    'E5' = '4.22%'
    'D6' = '2.096'
    'E6' = '-0.19%'
    'D7' = '0.9910'
    'E7' = '3.90%'
    'D8' = '2.564'
    'E8' = '-0.06%'
    'D9' = '0.1195'
    'E9' = '4.33%'
    'D10' = '0.1931'
    'E10' = '1.35%'
    'E11' = '1.54%'
    'D12' = '0.09916'
    'E12' = '-0.80%'
    'D13' = '0.04675'
    'E13' = '-3.35%'
    'E14' = '-0.54%'
    'D15' = '0.001289'
    'E15' = '0.72%'
    'D16' = '0.005931'
    'E16' = '-1.06%'
    'D17' = '3.394'
    'E17' = '0.74%'
    'D18' = '4.430'
    'E18' = '0.68%'
    'D19' = '0.3370'
    'E19' = '-0.65%'
    'D20' = '0.1354'
    'E20' = '-1.83%'
    'D21' = '0.2566'
    'E21' = '-0.65%'
    'D22' = '0.04133'
    'E22' = '1.24%'
    'D23' = '0.001293'
    'E23' = '1.55%'
    'D24' = '0.004541'
    'E24' = '4.31%'
    'E25' = '8.47%'
    'E26' = '0.02%'
    'D38' = '0.02698'
    'E38' = '4.24%'
    'D39' = '0.05755'
    'E39' = '-0.52%'
    'D40' = '0.007905'
    'E40' = '4.62%'
    'D41' = '0.1433'
    'E41' = '2.13%'
    'D42' = '0.007941'
    'E42' = '8.21%'
    'D43' = '0.002023'
    'E43' = '0.38%'
    'D44' = '0.008941'
    'E44' = '-1.53%'
    'D45' = '0.3409'
    'D46' = '0.00007041'
    'E46' = '0.63%'
    'E47' = '0.16%'
    'E48' = '0.29%'
    'D49' = '0.003536'
    'E49' = '0.16%'
    'D50' = '0.003391'
    'E50' = '-3.10%'
    'E51' = '0.16%'
}

foreach ($addr in $updates.Keys) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $updates[$addr]
}
